# Updates the worksheet date and every arithmetic problem in the table
# from the 2025-10-19 set to the 2025-10-20 set, one whole-text
# find-and-replace per run of text.
$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-10-19 Sunday", $true, $true, $false, $false, $false, $true, 1, $false, "2025-10-20 Monday", 2) | Out-Null
$d.Content.Find.Execute("22-13=", $true, $true, $false, $false, $false, $true, 1, $false, "81-37=", 2) | Out-Null
$d.Content.Find.Execute("34-11=", $true, $true, $false, $false, $false, $true, 1, $false, "30+58=", 2) | Out-Null
$d.Content.Find.Execute("78+18=", $true, $true, $false, $false, $false, $true, 1, $false, "48-22=", 2) | Out-Null
$d.Content.Find.Execute("73-16=", $true, $true, $false, $false, $false, $true, 1, $false, "14+71=", 2) | Out-Null
$d.Content.Find.Execute("55+21=", $true, $true, $false, $false, $false, $true, 1, $false, "29-6=", 2) | Out-Null
$d.Content.Find.Execute("44-17=", $true, $true, $false, $false, $false, $true, 1, $false, "65-26=", 2) | Out-Null
$d.Content.Find.Execute("25-8=", $true, $true, $false, $false, $false, $true, 1, $false, "21+64=", 2) | Out-Null
$d.Content.Find.Execute("7+61=", $true, $true, $false, $false, $false, $true, 1, $false, "92-8=", 2) | Out-Null
$d.Content.Find.Execute("59+36=", $true, $true, $false, $false, $false, $true, 1, $false, "18+22=", 2) | Out-Null
$d.Content.Find.Execute("65-2=", $true, $true, $false, $false, $false, $true, 1, $false, "23+32=", 2) | Out-Null
$d.Content.Find.Execute("16+33=", $true, $true, $false, $false, $false, $true, 1, $false, "67-34=", 2) | Out-Null
$d.Content.Find.Execute("33-31=", $true, $true, $false, $false, $false, $true, 1, $false, "83+15=", 2) | Out-Null
$d.Content.Find.Execute("48+9=", $true, $true, $false, $false, $false, $true, 1, $false, "15+21=", 2) | Out-Null
$d.Content.Find.Execute("60+30=", $true, $true, $false, $false, $false, $true, 1, $false, "2+40=", 2) | Out-Null
$d.Content.Find.Execute("76+7=", $true, $true, $false, $false, $false, $true, 1, $false, "90-44=", 2) | Out-Null
$d.Content.Find.Execute("38-38=", $true, $true, $false, $false, $false, $true, 1, $false, "23+54=", 2) | Out-Null
$d.Content.Find.Execute("22+33=", $true, $true, $false, $false, $false, $true, 1, $false, "11+63=", 2) | Out-Null
$d.Content.Find.Execute("26+56=", $true, $true, $false, $false, $false, $true, 1, $false, "26-3=", 2) | Out-Null
$d.Content.Find.Execute("91-0=", $true, $true, $false, $false, $false, $true, 1, $false, "77-59=", 2) | Out-Null
$d.Content.Find.Execute("5+85=", $true, $true, $false, $false, $false, $true, 1, $false, "45-20=", 2) | Out-Null
$d.Content.Find.Execute("41+40=", $true, $true, $false, $false, $false, $true, 1, $false, "4+5=", 2) | Out-Null
$d.Content.Find.Execute("24-2=", $true, $true, $false, $false, $false, $true, 1, $false, "32+44=", 2) | Out-Null
$d.Content.Find.Execute("33+8=", $true, $true, $false, $false, $false, $true, 1, $false, "98-89=", 2) | Out-Null
$d.Content.Find.Execute("36-12=", $true, $true, $false, $false, $false, $true, 1, $false, "64-20=", 2) | Out-Null
$d.Content.Find.Execute("45-21=", $true, $true, $false, $false, $false, $true, 1, $false, "73-36=", 2) | Out-Null
$d.Content.Find.Execute("94-0=", $true, $true, $false, $false, $false, $true, 1, $false, "70-36=", 2) | Out-Null
$d.Content.Find.Execute("93-59=", $true, $true, $false, $false, $false, $true, 1, $false, "0+98=", 2) | Out-Null
$d.Content.Find.Execute("92-83=", $true, $true, $false, $false, $false, $true, 1, $false, "32+56=", 2) | Out-Null
$d.Content.Find.Execute("0+70=", $true, $true, $false, $false, $false, $true, 1, $false, "79+18=", 2) | Out-Null
$d.Content.Find.Execute("87+0=", $true, $true, $false, $false, $false, $true, 1, $false, "69+7=", 2) | Out-Null
$d.Content.Find.Execute("45+33=", $true, $true, $false, $false, $false, $true, 1, $false, "49+31=", 2) | Out-Null
$d.Content.Find.Execute("3+23=", $true, $true, $false, $false, $false, $true, 1, $false, "82-20=", 2) | Out-Null
$d.Content.Find.Execute("11+26=", $true, $true, $false, $false, $false, $true, 1, $false, "83-40=", 2) | Out-Null
$d.Content.Find.Execute("18+59=", $true, $true, $false, $false, $false, $true, 1, $false, "49-29=", 2) | Out-Null
$d.Content.Find.Execute("95-20=", $true, $true, $false, $false, $false, $true, 1, $false, "24-11=", 2) | Out-Null
$d.Content.Find.Execute("3+7=", $true, $true, $false, $false, $false, $true, 1, $false, "72+15=", 2) | Out-Null
$d.Content.Find.Execute("63-54=", $true, $true, $false, $false, $false, $true, 1, $false, "75+23=", 2) | Out-Null
$d.Content.Find.Execute("44-34=", $true, $true, $false, $false, $false, $true, 1, $false, "88-3=", 2) | Out-Null
$d.Content.Find.Execute("20+3=", $true, $true, $false, $false, $false, $true, 1, $false, "44+4=", 2) | Out-Null
$d.Content.Find.Execute("25+48=", $true, $true, $false, $false, $false, $true, 1, $false, "14+69=", 2) | Out-Null
$d.Content.Find.Execute("26+53=", $true, $true, $false, $false, $false, $true, 1, $false, "30+54=", 2) | Out-Null
$d.Content.Find.Execute("72-53=", $true, $true, $false, $false, $false, $true, 1, $false, "73-21=", 2) | Out-Null
$d.Content.Find.Execute("16+55=", $true, $true, $false, $false, $false, $true, 1, $false, "59-3=", 2) | Out-Null
$d.Content.Find.Execute("15+19=", $true, $true, $false, $false, $false, $true, 1, $false, "71-5=", 2) | Out-Null
$d.Content.Find.Execute("45+45=", $true, $true, $false, $false, $false, $true, 1, $false, "41-32=", 2) | Out-Null
$d.Content.Find.Execute("9+79=", $true, $true, $false, $false, $false, $true, 1, $false, "3+33=", 2) | Out-Null
$d.Content.Find.Execute("94+1=", $true, $true, $false, $false, $false, $true, 1, $false, "0+67=", 2) | Out-Null
$d.Content.Find.Execute("41+4=", $true, $true, $false, $false, $false, $true, 1, $false, "26+25=", 2) | Out-Null
$d.Content.Find.Execute("97-46=", $true, $true, $false, $false, $false, $true, 1, $false, "68-40=", 2) | Out-Null
$d.Content.Find.Execute("81-68=", $true, $true, $false, $false, $false, $true, 1, $false, "61-33=", 2) | Out-Null
$d.Content.Find.Execute("28+48=", $true, $true, $false, $false, $false, $true, 1, $false, "49-26=", 2) | Out-Null
$d.Content.Find.Execute("27+53=", $true, $true, $false, $false, $false, $true, 1, $false, "27-9=", 2) | Out-Null
$d.Content.Find.Execute("68+31=", $true, $true, $false, $false, $false, $true, 1, $false, "1+37=", 2) | Out-Null
$d.Content.Find.Execute("32+64=", $true, $true, $false, $false, $false, $true, 1, $false, "92+2=", 2) | Out-Null
$d.Content.Find.Execute("51+13=", $true, $true, $false, $false, $false, $true, 1, $false, "94-24=", 2) | Out-Null
$d.Content.Find.Execute("3+53=", $true, $true, $false, $false, $false, $true, 1, $false, "91-39=", 2) | Out-Null
$d.Content.Find.Execute("57-6=", $true, $true, $false, $false, $false, $true, 1, $false, "95-46=", 2) | Out-Null
$d.Content.Find.Execute("72-44=", $true, $true, $false, $false, $false, $true, 1, $false, "33-28=", 2) | Out-Null
$d.Content.Find.Execute("39-3=", $true, $true, $false, $false, $false, $true, 1, $false, "94-65=", 2) | Out-Null
$d.Content.Find.Execute("61-50=", $true, $true, $false, $false, $false, $true, 1, $false, "57+41=", 2) | Out-Null
$d.Content.Find.Execute("22+42=", $true, $true, $false, $false, $false, $true, 1, $false, "57-50=", 2) | Out-Null
$d.Content.Find.Execute("24+39=", $true, $true, $false, $false, $false, $true, 1, $false, "52+27=", 2) | Out-Null
$d.Content.Find.Execute("12+19=", $true, $true, $false, $false, $false, $true, 1, $false, "39+44=", 2) | Out-Null
$d.Content.Find.Execute("53+9=", $true, $true, $false, $false, $false, $true, 1, $false, "73+6=", 2) | Out-Null
$d.Content.Find.Execute("79-22=", $true, $true, $false, $false, $false, $true, 1, $false, "54-33=", 2) | Out-Null
$d.Content.Find.Execute("98-58=", $true, $true, $false, $false, $false, $true, 1, $false, "7+0=", 2) | Out-Null
$d.Content.Find.Execute("63-60=", $true, $true, $false, $false, $false, $true, 1, $false, "45-15=", 2) | Out-Null
$d.Content.Find.Execute("24-13=", $true, $true, $false, $false, $false, $true, 1, $false, "83-65=", 2) | Out-Null
$d.Content.Find.Execute("78-41=", $true, $true, $false, $false, $false, $true, 1, $false, "98-23=", 2) | Out-Null
$d.Content.Find.Execute("54-12=", $true, $true, $false, $false, $false, $true, 1, $false, "69+8=", 2) | Out-Null
$d.Content.Find.Execute("85-48=", $true, $true, $false, $false, $false, $true, 1, $false, "8+62=", 2) | Out-Null
$d.Content.Find.Execute("29+21=", $true, $true, $false, $false, $false, $true, 1, $false, "82-27=", 2) | Out-Null
$d.Content.Find.Execute("50-39=", $true, $true, $false, $false, $false, $true, 1, $false, "2+82=", 2) | Out-Null
$d.Content.Find.Execute("44+35=", $true, $true, $false, $false, $false, $true, 1, $false, "98-84=", 2) | Out-Null
$d.Content.Find.Execute("53+4=", $true, $true, $false, $false, $false, $true, 1, $false, "42-3=", 2) | Out-Null
$d.Content.Find.Execute("8+18=", $true, $true, $false, $false, $false, $true, 1, $false, "35-33=", 2) | Out-Null
$d.Content.Find.Execute("84-14=", $true, $true, $false, $false, $false, $true, 1, $false, "72-60=", 2) | Out-Null
$d.Content.Find.Execute("20+7=", $true, $true, $false, $false, $false, $true, 1, $false, "31+2=", 2) | Out-Null
$d.Content.Find.Execute("52-46=", $true, $true, $false, $false, $false, $true, 1, $false, "74-14=", 2) | Out-Null
$d.Content.Find.Execute("70+4=", $true, $true, $false, $false, $false, $true, 1, $false, "6+58=", 2) | Out-Null
$d.Content.Find.Execute("97-36=", $true, $true, $false, $false, $false, $true, 1, $false, "50-9=", 2) | Out-Null
$d.Content.Find.Execute("50+1=", $true, $true, $false, $false, $false, $true, 1, $false, "18+14=", 2) | Out-Null
$d.Content.Find.Execute("59+37=", $true, $true, $false, $false, $false, $true, 1, $false, "25+29=", 2) | Out-Null
$d.Content.Find.Execute("20+39=", $true, $true, $false, $false, $false, $true, 1, $false, "84-10=", 2) | Out-Null
$d.Content.Find.Execute("98-14=", $true, $true, $false, $false, $false, $true, 1, $false, "42-10=", 2) | Out-Null
$d.Content.Find.Execute("95-38=", $true, $true, $false, $false, $false, $true, 1, $false, "57+28=", 2) | Out-Null
$d.Content.Find.Execute("97-30=", $true, $true, $false, $false, $false, $true, 1, $false, "4+15=", 2) | Out-Null
$d.Content.Find.Execute("38+43=", $true, $true, $false, $false, $false, $true, 1, $false, "79-61=", 2) | Out-Null
$d.Content.Find.Execute("58+5=", $true, $true, $false, $false, $false, $true, 1, $false, "83-82=", 2) | Out-Null
$d.Content.Find.Execute("78-76=", $true, $true, $false, $false, $false, $true, 1, $false, "5+76=", 2) | Out-Null
$d.Content.Find.Execute("48-17=", $true, $true, $false, $false, $false, $true, 1, $false, "71-51=", 2) | Out-Null
$d.Content.Find.Execute("79-64=", $true, $true, $false, $false, $false, $true, 1, $false, "7+65=", 2) | Out-Null
$d.Content.Find.Execute("1+92=", $true, $true, $false, $false, $false, $true, 1, $false, "89-30=", 2) | Out-Null
$d.Content.Find.Execute("6+70=", $true, $true, $false, $false, $false, $true, 1, $false, "67+12=", 2) | Out-Null
$d.Content.Find.Execute("31-28=", $true, $true, $false, $false, $false, $true, 1, $false, "81-28=", 2) | Out-Null
$d.Content.Find.Execute("65-19=", $true, $true, $false, $false, $false, $true, 1, $false, "5-0=", 2) | Out-Null
$d.Content.Find.Execute("83-83=", $true, $true, $false, $false, $false, $true, 1, $false, "17+25=", 2) | Out-Null
$d.Content.Find.Execute("44-23=", $true, $true, $false, $false, $false, $true, 1, $false, "42+40=", 2) | Out-Null
$d.Content.Find.Execute("85-20=", $true, $true, $false, $false, $false, $true, 1, $false, "42+31=", 2) | Out-Null
$d.Content.Find.Execute("49-7=", $true, $true, $false, $false, $false, $true, 1, $false, "60-17=", 2) | Out-Null
